$d = $word.ActiveDocument

# Common package wrapper used for all InsertXML calls (Range.InsertXML expects a
# full WordOpenXML single-file package payload).
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-ParagraphXml($paragraphIndex, $innerBodyXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    $xml = $pkgOpen + $innerBodyXml + $pkgClose
    $r.InsertXML($xml)
}

$rpr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# ---------------------------------------------------------------------------
# Paragraph 3: "...We spent about 35 hours total..." -> split "hours" out with
# gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$p3 = '<w:p w:rsidR="00653DCB" w:rsidRDefault="00534CE1" w:rsidP="004F1554">' + `
  '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' + $rpr + '</w:pPr>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve">The customer asked us to implement </w:t></w:r>' + `
  '<w:r w:rsidR="007534EC">' + $rpr + '<w:t xml:space="preserve">the </w:t></w:r>' + `
  '<w:r w:rsidR="007534EC" w:rsidRPr="00B30383"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>account creation</w:t></w:r>' + `
  '<w:r w:rsidR="007534EC">' + $rpr + '<w:t xml:space="preserve"> and </w:t></w:r>' + `
  '<w:r w:rsidR="007534EC" w:rsidRPr="00B30383"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>single-player</w:t></w:r>' + `
  '<w:r w:rsidR="007534EC">' + $rpr + '<w:t xml:space="preserve"> user stories.</w:t></w:r>' + `
  '<w:r w:rsidR="00925110">' + $rpr + '<w:t xml:space="preserve">  We implemented these stories, as well </w:t></w:r>' + `
  '<w:r w:rsidR="00F926ED">' + $rpr + '<w:t>as the user-submitted questions.  W</w:t></w:r>' + `
  '<w:r w:rsidR="00890194">' + $rpr + '<w:t xml:space="preserve">e spent about 35 </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r>' + $rpr + '<w:t>hours</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidR="00CB6B9B">' + $rpr + '<w:t>total working on the project; compared to our estimate, we were spot on.</w:t></w:r>' + `
  '</w:p>'

Replace-ParagraphXml 3 $p3

# ---------------------------------------------------------------------------
# Paragraph 9: "Pair Programming" bullet -> split "php" out with
# spellStart/spellEnd proofErr markers.
# ---------------------------------------------------------------------------
$p9 = '<w:p w:rsidR="00BC7539" w:rsidRDefault="00BC7539" w:rsidP="001043E6">' + `
  '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>' + $rpr + '</w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Pair Programming</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidR="00893A27">' + $rpr + '<w:t>–</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidR="00893A27">' + $rpr + '<w:t xml:space="preserve">we split into pairs, with one person of the pair coding, the other checking.  </w:t></w:r>' + `
  '<w:r w:rsidR="007E54DF">' + $rpr + '<w:t xml:space="preserve">Then the other would code and the other would error check.  </w:t></w:r>' + `
  '<w:r w:rsidR="00893A27">' + $rpr + '<w:t>One pair worke</w:t></w:r>' + `
  '<w:r w:rsidR="007E54DF">' + $rpr + '<w:t xml:space="preserve">d on the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + $rpr + '<w:t>php</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve">, the other on the JavaScript.  Then the pairs swapped assignments.  We went from there.</w:t></w:r>' + `
  '</w:p>'

Replace-ParagraphXml 9 $p9

# ---------------------------------------------------------------------------
# Paragraph 11: "Sustainable Pace" bullet -> drop the trailing _GoBack bookmark
# (it moves to the very end of the new paragraph added after item 14).
# ---------------------------------------------------------------------------
$p11 = '<w:p w:rsidR="00011240" w:rsidRDefault="00011240" w:rsidP="001043E6">' + `
  '<w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>' + $rpr + '</w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>Sustainable Pace</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidR="00886130">' + $rpr + '<w:t>–</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r w:rsidR="00886130">' + $rpr + '<w:t>our group did not work past 10 PM each day of work.  And we did not exceed 8 hours of work per each workday.</w:t></w:r>' + `
  '</w:p>'

Replace-ParagraphXml 11 $p11

# ---------------------------------------------------------------------------
# Paragraph 14: single-player test-case bullet is unchanged in content, but we
# must append a brand-new ListParagraph bullet right after it (describing the
# user-submitted-question test case), carrying the relocated _GoBack bookmark
# at its very end.
# ---------------------------------------------------------------------------
$p14 = '<w:p w:rsidR="0061347A" w:rsidRPr="0061347A" w:rsidRDefault="00087BB8" w:rsidP="00FC6F03">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/>' + $rpr + '</w:pPr>' + `
  '<w:r>' + $rpr + '<w:t>For the single-player story, we tested that the trivia game could be played; that is, a random question is presented to the player with four choices to pick as the correct answer.  After submitting their answer, an alert pops up and tells them whether their answer is correct or not.</w:t></w:r>' + `
  '<w:r w:rsidR="007E6C65">' + $rpr + '<w:t xml:space="preserve">  Then a new question is presented and the cycle continues.</w:t></w:r>' + `
  '<w:r w:rsidR="00EC40FE">' + $rpr + '<w:t xml:space="preserve">  The input is the players answer choice.  The expected output is the alert box indicating right/wrong.</w:t></w:r>' + `
  '</w:p>'

$p15 = '<w:p>' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/>' + $rpr + '</w:pPr>' + `
  '<w:r>' + $rpr + '<w:t>The user-subm</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t>itted questions were tested for whether or not they could be stored correctly in the database and then generated for the player to answer during single-player mode.</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve">  The inputs are</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t>&#8220;Question&#8221;, &#8220;Correct Answer&#8221;, and three &#8220;Fake Answer&#8221; fields.</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve">  </w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:lastRenderedPageBreak/><w:t>The expected output is &#8220;Congratulations!  Your question was successfully submitted!&#8221;</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> that is displayed</w:t></w:r>' + `
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> at the top of the page.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'

$p14and15 = $p14 + $p15
Replace-ParagraphXml 14 $p14and15

Write-Output "done"
